# Update the "想去人数" (column F) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Row -> new F value mapping (same update applies to both sheets)
$updates = @{
    3  = 1711
    6  = 34
    7  = 11947
    8  = 42
    11 = 409
    13 = 852
    14 = 13467
    15 = 13443
    16 = 39
    17 = 154
    18 = 20
    19 = 38
    20 = 288
    23 = 147
    24 = 170
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
